$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.955.13'
$ws.Range('E2').Value = '  -0.71%  '
$ws.Range('D3').Value = '1.646.30'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('D5').Value = "'217.91"
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').Value = "'0.5240"
$ws.Range('E6').Value = '  +0.69%  '
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('D8').Value = "'0.2617"
$ws.Range('E8').Value = '  -1.91%  '
$ws.Range('D9').Value = "'0.06278"
$ws.Range('E9').Value = '  -0.77%  '
$ws.Range('D10').Value = "'20.33"
$ws.Range('E10').Value = '  -3.58%  '
$ws.Range('D11').Value = "'0.07734"
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = "'4.457"
$ws.Range('E12').Value = '  +0.36%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.567.65'
$ws.Range('E13').Value = '  -5.25%  '
$ws.Range('D14').Value = "'0.5442"
$ws.Range('E14').Value = '  -0.26%  '
$ws.Range('D15').Value = '0.0₅8081'
$ws.Range('E15').Value = '  -1.74%  '
$ws.Range('D16').Value = "'64.89"
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('D17').Value = '25.987.92'
$ws.Range('E17').Value = '  -0.80%  '
$ws.Range('E19').Value = '  -2.60%  '
$ws.Range('D20').Value = "'191.85"
$ws.Range('E20').Value = '  -0.41%  '
$ws.Range('D21').Value = "'10.03"
$ws.Range('E21').Value = '  -1.22%  '
$ws.Range('D22').Value = "'5.971"
$ws.Range('E22').Value = '  -2.07%  '
$ws.Range('D23').Value = "'1.004"
$ws.Range('E23').Value = '  -0.32%  '
$ws.Range('D24').Value = "'139.47"
$ws.Range('E24').Value = '  +1.48%  '
$ws.Range('D25').Value = "'0.1237"
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').Value = "'7.261"
$ws.Range('E26').Value = '  +0.43%  '
$ws.Range('E27').Value = '  +0.36%  '
$ws.Range('D28').Value = "'1.422"
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('D29').Value = "'0.05935"
$ws.Range('E29').Value = '  -1.74%  '
$ws.Range('D30').Value = "'1.274"
$ws.Range('E30').Value = '  -0.79%  '
$ws.Range('D31').Value = "'3.487"
$ws.Range('E31').Value = '  -2.12%  '
$ws.Range('D32').Value = "'3.235"
$ws.Range('E32').Value = '  -3.14%  '
$ws.Range('D33').Value = "'1.528"
$ws.Range('E33').Value = '  -7.39%  '
$ws.Range('D34').Value = "'2.412"
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('D35').Value = "'0.9394"
$ws.Range('E35').Value = '  -4.26%  '
$ws.Range('D36').Value = "'2.743"
$ws.Range('E36').Value = '  -1.31%  '
$ws.Range('D37').Value = "'0.5722"
$ws.Range('E37').Value = '  -3.23%  '
$ws.Range('D38').Value = "'0.01602"
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('D39').Value = "'5.860"
$ws.Range('E39').Value = '  -1.62%  '
$ws.Range('D40').Value = "'0.8445"
$ws.Range('E40').Value = '  -2.24%  '
$ws.Range('D42').Value = "'100.62"
$ws.Range('E42').Value = '  +0.86%  '
$ws.Range('D43').Value = '1.004.27'
$ws.Range('E43').Value = '  -3.03%  '
$ws.Range('D44').Value = '1.785.84'
$ws.Range('E44').Value = '  -0.59%  '
$ws.Range('D45').Value = "'56.63"
$ws.Range('E45').Value = '  -0.92%  '
$ws.Range('E46').Value = '  -2.49%  '
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('E48').Value = '  +1.29%  '
$ws.Range('D49').Value = "'1.477"
$ws.Range('E49').Value = '  +0.82%  '
$ws.Range('D50').Value = "'0.05147"
$ws.Range('E50').Value = '  -0.55%  '
$ws.Range('D51').Value = "'7.822"
$ws.Range('E51').Value = '  -3.32%  '
